# Update the Markov transition-probability matrix on Sheet1 (A1:S19).
# These new probabilities reflect the expanded game sample used to
# recompute the team-specific transition matrix (more games simulated,
# faster simulate-game logic).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2259259259259259
$ws.Range("C2").Value = 0.5074074074074074
$ws.Range("J2").Value = 0.01481481481481482
$ws.Range("P2").Value = 0.1407407407407407
$ws.Range("S2").Value = 0.1111111111111111
$ws.Range("B3").Value = 0.006622516556291391
$ws.Range("C3").Value = 0.03973509933774835
$ws.Range("J3").Value = 0.05298013245033113
$ws.Range("P3").Value = 0.7417218543046358
$ws.Range("S3").Value = 0.1589403973509934
$ws.Range("J4").Value = 0.09090909090909091
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.1590909090909091
$ws.Range("B6").Value = 0.04878048780487805
$ws.Range("D6").Value = 0.01463414634146342
$ws.Range("F6").Value = 0.03902439024390244
$ws.Range("J6").Value = 0.2731707317073171
$ws.Range("O6").Value = 0.03414634146341464
$ws.Range("Q6").Value = 0.1414634146341463
$ws.Range("R6").Value = 0.08292682926829269
$ws.Range("S6").Value = 0.3658536585365854
$ws.Range("B7").Value = 0.1215469613259668
$ws.Range("D7").Value = 0.02762430939226519
$ws.Range("F7").Value = 0.03867403314917127
$ws.Range("J7").Value = 0.143646408839779
$ws.Range("O7").Value = 0.005524861878453038
$ws.Range("Q7").Value = 0.2430939226519337
$ws.Range("R7").Value = 0.1104972375690608
$ws.Range("S7").Value = 0.3093922651933702
$ws.Range("B8").Value = 0.1184210526315789
$ws.Range("D8").Value = 0.02368421052631579
$ws.Range("F8").Value = 0.05263157894736842
$ws.Range("J8").Value = 0.1052631578947368
$ws.Range("O8").Value = 0.01052631578947368
$ws.Range("Q8").Value = 0.2131578947368421
$ws.Range("R8").Value = 0.09210526315789473
$ws.Range("S8").Value = 0.3842105263157894
$ws.Range("B9").Value = 0.103960396039604
$ws.Range("D9").Value = 0.009900990099009901
$ws.Range("F9").Value = 0.06930693069306931
$ws.Range("J9").Value = 0.1534653465346535
$ws.Range("O9").Value = 0.01485148514851485
$ws.Range("Q9").Value = 0.2227722772277228
$ws.Range("R9").Value = 0.0891089108910891
$ws.Range("S9").Value = 0.3366336633663367
$ws.Range("B10").Value = 0.08498349834983498
$ws.Range("D10").Value = 0.02062706270627063
$ws.Range("F10").Value = 0.07508250825082509
$ws.Range("J10").Value = 0.146039603960396
$ws.Range("O10").Value = 0.01072607260726073
$ws.Range("Q10").Value = 0.2409240924092409
$ws.Range("R10").Value = 0.08663366336633663
$ws.Range("S10").Value = 0.334983498349835
$ws.Range("G11").Value = 0.1342281879194631
$ws.Range("J11").Value = 0.1140939597315436
$ws.Range("K11").Value = 0.1912751677852349
$ws.Range("L11").Value = 0.5536912751677853
$ws.Range("S11").Value = 0.006711409395973154
$ws.Range("G12").Value = 0.7093023255813954
$ws.Range("J12").Value = 0.2151162790697674
$ws.Range("K12").Value = 0.01744186046511628
$ws.Range("L12").Value = 0.04069767441860465
$ws.Range("S12").Value = 0.01744186046511628
$ws.Range("G13").Value = 0.8709677419354839
$ws.Range("J13").Value = 0.1290322580645161
$ws.Range("F15").Value = 0.04166666666666666
$ws.Range("H15").Value = 0.1614583333333333
$ws.Range("I15").Value = 0.09375
$ws.Range("J15").Value = 0.359375
$ws.Range("K15").Value = 0.03645833333333334
$ws.Range("M15").Value = 0.005208333333333333
$ws.Range("N15").Value = 0.01041666666666667
$ws.Range("O15").Value = 0.046875
$ws.Range("S15").Value = 0.2447916666666667
$ws.Range("F16").Value = 0.01685393258426966
$ws.Range("H16").Value = 0.1573033707865168
$ws.Range("I16").Value = 0.1123595505617977
$ws.Range("J16").Value = 0.3146067415730337
$ws.Range("K16").Value = 0.1348314606741573
$ws.Range("M16").Value = 0.005617977528089887
$ws.Range("N16").Value = 0.005617977528089887
$ws.Range("O16").Value = 0.07303370786516854
$ws.Range("S16").Value = 0.1797752808988764
$ws.Range("F17").Value = 0.01646090534979424
$ws.Range("H17").Value = 0.154320987654321
$ws.Range("I17").Value = 0.09259259259259259
$ws.Range("J17").Value = 0.4176954732510288
$ws.Range("K17").Value = 0.1069958847736626
$ws.Range("M17").Value = 0.01440329218106996
$ws.Range("O17").Value = 0.05761316872427984
$ws.Range("S17").Value = 0.139917695473251
$ws.Range("F18").Value = 0.01020408163265306
$ws.Range("H18").Value = 0.1581632653061225
$ws.Range("I18").Value = 0.0663265306122449
$ws.Range("J18").Value = 0.413265306122449
$ws.Range("K18").Value = 0.1326530612244898
$ws.Range("M18").Value = 0.03061224489795918
$ws.Range("O18").Value = 0.05612244897959184
$ws.Range("S18").Value = 0.1326530612244898
$ws.Range("F19").Value = 0.01876675603217158
$ws.Range("H19").Value = 0.1885612153708668
$ws.Range("I19").Value = 0.09204647006255585
$ws.Range("J19").Value = 0.3663985701519213
$ws.Range("K19").Value = 0.1117068811438785
$ws.Range("M19").Value = 0.0160857908847185
$ws.Range("N19").Value = 0.002680965147453083
$ws.Range("O19").Value = 0.07327971403038427
$ws.Range("S19").Value = 0.13047363717605

Write-Host "Applied 108 cell updates"